$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70, pushing the existing rows 70-105 down to 71-106.
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new weekly price observation.
$ws.Range("A70").Value = 1
$ws.Range("B70").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C70").Value = "Arica y Parinacota"
$ws.Range("D70").Value = 44846
$ws.Range("E70").Value = 15
$ws.Range("F70").Value = "Fruta"
$ws.Range("G70").Value = 100102
$ws.Range("H70").Value = "Cítricos"
$ws.Range("I70").Value = 100102005
$ws.Range("J70").Value = "Naranja"
$ws.Range("K70").Value = "Lane Late"
$ws.Range("L70").Value = "Segunda"
$ws.Range("M70").Value = 270
$ws.Range("N70").Value = 700
$ws.Range("O70").Value = 750
$ws.Range("P70").Value = 725
$ws.Range("Q70").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R70").Value = "Región de O'Higgins"
$ws.Range("S70").Value = 725
$ws.Range("T70").Value = 1
